# Update capital structure / earnings database values for the Thailand
# reinsurance sheet (rows 2-4 of "earnings_debt"), per the refreshed source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Thaire Insurance Public Company Limited
$ws.Range("D2").Value = -0.03335
$ws.Range("E2").Value = -0.3295
$ws.Range("G2").Value = 0.06637085474077534
$ws.Range("H2").Value = 0.06637085474077534
$ws.Range("I2").Value = 0.04742816282245576
$ws.Range("J2").Value = 0.04016876792942909
$ws.Range("K2").Value = 8.43
$ws.Range("L2").Value = 0.03937412424100887
$ws.Range("M2").Value = 2.47
$ws.Range("N2").Value = 0.01016879374228078
$ws.Range("O2").Value = 0.2930011862396205
$ws.Range("P2").Value = 2.47
$ws.Range("Q2").Value = 0.01016879374228078
$ws.Range("R2").Value = 0.2930011862396205
$ws.Range("U2").Value = 11.78
$ws.Range("V2").Value = 0.04849732400164677
$ws.Range("W2").Value = 0.0521881526323692
$ws.Range("X2").Value = 0.08787038733482352
$ws.Range("Y2").Value = -0.03568223470245432
$ws.Range("Z2").Value = 1.329110707861179
$ws.Range("AA2").Value = 0.05488373020113378
$ws.Range("AB2").Value = 0.08776888700327277
$ws.Range("AC2").Value = -0.03288515680213899
$ws.Range("AD2").Value = 0.255
$ws.Range("AE2").Value = 0.1331516985611
$ws.Range("AF2").Value = 0.3881516985611
$ws.Range("AG2").Value = -11.3918483014389
$ws.Range("AH2").Value = 0.001595440204757804
$ws.Range("AI2").Value = 0.00252393760035502
$ws.Range("AJ2").Value = -0.0492071152478114
$ws.Range("AK2").Value = -0.08021967869577117
$ws.Range("AL2").Value = 0.015
$ws.Range("AM2").Value = 0.015
$ws.Range("AN2").Value = 0.02303314967030982
$ws.Range("AO2").Value = 674.6666666666667
$ws.Range("AP2").Value = -1.028980968425517
$ws.Range("AQ2").Value = 674.6666666666667

# Row 3 - Thaire Life Assurance Public Company Limited (SET:THREL)
$ws.Range("D3").Value = 0.0413
$ws.Range("E3").Value = -0.261
$ws.Range("G3").Value = 0.1077235772357724
$ws.Range("H3").Value = 0.1077235772357724
$ws.Range("I3").Value = 0.04352973951119445
$ws.Range("J3").Value = 0.03833073300531198
$ws.Range("K3").Value = 2.81
$ws.Range("L3").Value = 0.03807588075880759
$ws.Range("M3").Value = 2.47
$ws.Range("N3").Value = 0.03908227848101266
$ws.Range("O3").Value = 0.8790035587188613
$ws.Range("P3").Value = 2.47
$ws.Range("Q3").Value = 0.03908227848101266
$ws.Range("R3").Value = 0.8790035587188613
$ws.Range("U3").Value = 1.18
$ws.Range("V3").Value = 0.01867088607594937
$ws.Range("W3").Value = 0.05758196721311476
$ws.Range("X3").Value = 0.08788685516694154
$ws.Range("Y3").Value = -0.03030488795382678
$ws.Range("Z3").Value = 1.573811566699595
$ws.Range("AA3").Value = 0.0603253509638339
$ws.Range("AB3").Value = 0.08776546749707485
$ws.Range("AC3").Value = -0.02744011653324095
$ws.Range("AD3").Value = 0.053
$ws.Range("AE3").Value = 0.07252612036925087
$ws.Range("AF3").Value = 0.1255261203692509
$ws.Range("AG3").Value = -1.054473879630749
$ws.Range("AH3").Value = 0.00198223572798512
$ws.Range("AI3").Value = 0.002994026121672563
$ws.Range("AJ3").Value = -0.01696781643763617
$ws.Range("AK3").Value = -0.02587950089331656
$ws.Range("AL3").Value = 0.003
$ws.Range("AM3").Value = 0.003
$ws.Range("AN3").Value = 0.01569440331655315
$ws.Range("AO3").Value = 1066.666666666667
$ws.Range("AP3").Value = -0.3122516670508585
$ws.Range("AQ3").Value = 1066.666666666667

# Row 4 - Thai Reinsurance Public Company Limited (SET:THRE)
$ws.Range("D4").Value = -0.108
$ws.Range("E4").Value = -0.398
$ws.Range("G4").Value = 0.04461867426942266
$ws.Range("H4").Value = 0.04461867426942266
$ws.Range("I4").Value = 0.04947879461412423
$ws.Range("J4").Value = 0.04024179822451204
$ws.Range("K4").Value = 5.62
$ws.Range("L4").Value = 0.04005702066999287
$ws.Range("O4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("U4").Value = 10.6
$ws.Range("V4").Value = 0.05898720089037285
$ws.Range("W4").Value = 0.04679433805162365
$ws.Range("X4").Value = 0.08785391950270549
$ws.Range("Y4").Value = -0.04105958145108184
$ws.Range("Z4").Value = 1.228625747850342
$ws.Range("AA4").Value = 0.04944210943843366
$ws.Range("AB4").Value = 0.0877723065094707
$ws.Range("AC4").Value = -0.03833019707103704
$ws.Range("AD4").Value = 0.202
$ws.Range("AE4").Value = 0.06062557819184914
$ws.Range("AF4").Value = 0.2626255781918492
$ws.Range("AG4").Value = -10.33737442180815
$ws.Range("AH4").Value = 0.001459333999757306
$ws.Range("AI4").Value = 0.00252393760035502
$ws.Range("AJ4").Value = -0.06103692822732935
$ws.Range("AK4").Value = -0.1020847954789198
$ws.Range("AL4").Value = 0.012
$ws.Range("AM4").Value = 0.012
$ws.Range("AN4").Value = 0.02625422407070445
$ws.Range("AO4").Value = 576.6666666666666
$ws.Range("AP4").Value = -1.028980968425517
$ws.Range("AQ4").Value = 576.6666666666666
